# Added program coding phase report
# Populates rows 10-12 of the "Coding Phase Defects" sheet with the
# first set of coding-phase review defects, and makes that sheet the
# active / selected tab (matching the author's last saved selection).

$wb = $excel.ActiveWorkbook

$codingSheet = $wb.Worksheets.Item("Coding Phase Defects")

# Row 10: C06 / File Repos / Input data is not obtained due to incorrect input files paths
$codingSheet.Range("C10").Value = "C06"
$codingSheet.Range("D10").Value = "File Repos"
$codingSheet.Range("E10").Value = "Input data is not obtained due to incorrect input files paths"
$codingSheet.Rows.Item(10).RowHeight = 30

# Row 11: C12 / Activity / Duration is a time ...
$codingSheet.Range("C11").Value = "C12"
$codingSheet.Range("D11").Value = "Activity"
$codingSheet.Range("E11").Value = "Duration is a time and it is declared as a date for e.g. but there are so many things to take into account here regarding variable naming, redundant null initializations and so on and so forth"
$codingSheet.Rows.Item(11).RowHeight = 75

# Row 12: C01 / Everywhere / Architectural decisions mainly ...
# (shared-string table records C, then the comment text, then "Everywhere")
$codingSheet.Range("C12").Value = "C01"
$codingSheet.Range("E12").Value = "Architectural decisions mainly, which consist of bad packaging and layering, wrong location for many classes and/or interfaces, for e.g. the repositories which are in a model package plus some weird output messages decisions and so on"
$codingSheet.Range("D12").Value = "Everywhere"
$codingSheet.Rows.Item(12).RowHeight = 90

# Make the Coding Phase Defects sheet the active tab and restore its
# last-used selection, just as it was when the author saved the file.
$codingSheet.Activate()
$codingSheet.Range("D13").Select()
